# Add two new norm-rule rows ("hyperlink3" and "hyperlink4") to the
# "Normative Rules" table, right after the existing "hyperlink2" row,
# so the norm-rule exporter's handling of <<foo>> / <<foo,link text>>
# anchors for non-norm targets is covered by a test fixture row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# "hyperlink2" lives in row 22 (row 1 is the header). Insert two blank
# worksheet rows right after it -- this pushes table1/table2/... etc.
# down by two rows and grows the sheet's used range accordingly.
$ws.Rows.Item(23).Insert()
$ws.Rows.Item(24).Insert()

# The ListObject (Table1) doesn't auto-grow from a plain row insert, so
# resize it explicitly to re-include the two new rows (old range was
# A1:F44 -> new range A1:F46, same as the sheet's "last row is outside
# the table" pattern is preserved for the trailing admon5/parameter row).
$lo.Resize($ws.Range("A1:F46"))

# Row 23: hyperlink3
$ws.Range("A23").Value = "my-chapter_name"
$ws.Range("B23").Value = "hyperlink3"
$ws.Range("C23").Value = "ABC &lt;&lt;non-norm-anchor&gt;&gt; DEF"
$ws.Range("D23").Value = '["norm:hyperlink3"]'

# Row 24: hyperlink4
$ws.Range("A24").Value = "my-chapter_name"
$ws.Range("B24").Value = "hyperlink4"
$ws.Range("C24").Value = "DEF &lt;&lt;non-norm-anchor,custom text&gt;&gt; GHI"
$ws.Range("D24").Value = '["norm:hyperlink4"]'
